# Insert a new September transaction row in the "2024" sheet.
# This pushes the existing rows 36-130 down to 37-131 (matching the
# growth of the sheet's used range from A1:Y130 to A1:Y131), and
# populates the newly inserted row 36 with the new transaction.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("2024")

# Insert a new blank row above row 36 - shifts rows 36:130 -> 37:131.
$ws.Rows.Item(36).Insert()

# Populate the new row's September Details / Date columns.
$ws.Range("R36").Value2 = "money google icici"
$ws.Range("S36").Value2 = "2024-09-10 20:42:12"
